$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# --- Simple value updates -------------------------------------------------

# Version (row 3): 0.2.1 -> 0.2.2
$ws.Range("B3").Value = "0.2.2"

# Date (row 8): refreshed build date
$ws.Range("B8").Value = "2024-09-11T16:17:59-05:00"

# Contact (row 10): now shows the publisher contact instead of the old
# "No display for ContactDetail" placeholder.
$ws.Range("B10").Value = "MITRE (https://www.mitre.org)"

# --- Insert a new "Jurisdiction" row after "Contact" ----------------------
# (pushes "Description" and everything below it down by one row, so the
# table grows from A1:B21 to A1:B22)

# First, give the about-to-exist row 22 the right formatting (border/
# alignment/font) by copying the format from the current last row (21).
$ws.Range("A21:B21").Copy() | Out-Null
$ws.Range("A22:B22").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# Stage a snapshot of rows 11:21 in a scratch area far away from both the
# copy source and the paste destination, so the upcoming overlapping
# (source rows 12:21 == destination rows 12:21) shift can't read back
# already-overwritten cells.
$ws.Range("A11:B21").Copy() | Out-Null
$ws.Range("Z1:AA11").PasteSpecial(-4104) | Out-Null  # xlPasteAll
$excel.CutCopyMode = $false

# Clear the destination block so that blank cells in the source actually
# blank out the destination (paste alone leaves old content behind where
# the incoming cell is empty).
$ws.Range("A12:B22").ClearContents() | Out-Null

# Paste the snapshot down into rows 12:22.
$ws.Range("Z1:AA11").Copy() | Out-Null
$ws.Range("A12:B22").PasteSpecial(-4104) | Out-Null  # xlPasteAll
$excel.CutCopyMode = $false

# Remove the scratch snapshot.
$ws.Range("Z1:AA11").Clear() | Out-Null

# Populate the newly freed row 11 with the new "Jurisdiction" property.
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
